$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New data rows for "Remessa" / "Material" / "Quantidade".
# Rows 53-54 held the two prior entries (84004823/15211-DLO-I and
# 84004824/10251-ARI-I); they are pushed down to rows 60-61, the
# previously-empty rows 55-63 are populated, and the newly-added
# shipments occupy the freed rows 53-59.
$data = @(
    @("80266554", "10246-ARI-I", 1),
    @("80266555", "10256-ARI-I", 1),
    @("80266556", "10377-ARI-I", 1),
    @("80266557", "10119-ATE-I", 5),
    @("80266557", "23359-ATE-I", 5),
    @("80266557", "23380-ATE-I", 2),
    @("80266563", "31497-MNW-I", 10),
    @("84004823", "15211-DLO-I", 1),
    @("84004824", "10251-ARI-I", 1),
    @("84004825", "15211-DLO-I", 2),
    @("84004825", "15386-DLO-I", 1)
)

# Column A ("Remessa") values are purely numeric strings (e.g. "80266554").
# Pre-format that range as Text so Excel keeps them as text instead of
# silently converting to numbers.
$ws.Range("A53:A63").NumberFormat = "@"

$row = 53
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

$ws.Range("F13").Select()
